# "Ran code for averaged intensities on spiral schemes"
#
# The averaged-intensities table (UniformA sheet) is regenerated with extra
# sampling schemes. Relative to the previous run:
#   - "Gaussian-Quadrature" (previously the LAST scheme row) is now listed
#     right after "Ring Perpendicular to TD" (i.e. it moved up in the table).
#   - Three brand new schemes follow it: "Spiral-90deg-10rot-5space",
#     "Spiral-90deg-15rot-5space", "Spiral-90deg-10rot-3space".
#   - Every other scheme (NoRotation-tilt60deg ... HexGrid-60degTilt5degRes)
#     keeps its previous relative order, just pushed down by the 4 rows
#     above.
#   - Column A is just the running HKL/scheme index (0, 1, 2, ...), and
#     every data cell (C:P) is 1 for every scheme row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Make room: insert 4 blank rows right before the current row 10
#    ("NoRotation-tilt60deg"). This pushes the existing rows 10-16 down to
#    14-20, including the old "Gaussian-Quadrature" row (16 -> 20).
$ws.Rows("10:13").Insert()

# 2) Move the old "Gaussian-Quadrature" row (now at row 20) up into the
#    first of the freshly inserted blank rows (row 10).
$ws.Rows("20").Cut()
$ws.Rows("10").PasteSpecial(-4104)  # xlPasteAll
$ws.Rows("20").Delete()

# 3) Re-apply the bold/bordered/centered style used by every column-A cell
#    (copy it from A2, which still carries the original style) across the
#    whole data range, since the row insert/cut above can leave the moved
#    / inserted rows with a slightly different style id.
$ws.Range("A2").Copy()
$ws.Range("A3:A19").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 4) Fill in the new Spiral-* rows (11, 12, 13) and the relocated
#    Gaussian-Quadrature row's label/index (row 10 already has the correct
#    data values pasted from the cut, just fix its HKL index below).
$schemeRows = @(
    @{ Row = 10; Name = "Gaussian-Quadrature" },
    @{ Row = 11; Name = "Spiral-90deg-10rot-5space" },
    @{ Row = 12; Name = "Spiral-90deg-15rot-5space" },
    @{ Row = 13; Name = "Spiral-90deg-10rot-3space" }
)

foreach ($row in $schemeRows) {
    $ws.Cells.Item($row.Row, 2).Value = $row.Name
    for ($col = 3; $col -le 16; $col++) {
        $ws.Cells.Item($row.Row, $col).Value = 1
    }
}

# 5) Renumber the HKL index column (A) sequentially for every data row
#    (A2 = 0 for the header row, A3:A19 = 1..17).
for ($r = 2; $r -le 19; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}
